$d = $word.ActiveDocument

# The document's Date paragraph (pStyle "Date") reads "October 18, 2017".
# Update it to "November 1, 2017" by replacing just the two words that
# changed, scoped to that paragraph so nothing else in the document is
# touched.
$datePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Date") {
        $datePara = $p
        break
    }
}

$rng = $datePara.Range
$rng.Find.ClearFormatting()
$rng.Find.Execute("October", $false, $true, $false, $false, $false, $true, 1, $false, "November", 2)

$rng = $datePara.Range
$rng.Find.ClearFormatting()
$rng.Find.Execute("18,", $false, $true, $false, $false, $false, $true, 1, $false, "1,", 2)
